# Applies the two edits described by the diff:
#
#  1. Remove the stray "_GoBack" bookmark that sat right after the
#     "...-серверах." run (it is not re-created at that spot).
#
#  2. Merge the "Дополнительные требования(...)" paragraph with the
#     following "6.1 Адаптивная верстка..." paragraph: the "Дополнительные
#     требования(...)" paragraph (and its "6.1 " numbering prefix) is
#     dropped, and the surviving paragraph (which keeps the original
#     "Дополнительные требования" paragraph's numbered-list formatting)
#     ends up with the single run text
#     " Адаптивная верстка страниц для мобильного приложения.", preceded
#     by a fresh "_GoBack" bookmark (the one removed in step 1 effectively
#     re-lands here, which is where Word leaves it after this kind of
#     in-place edit).

$d = $word.ActiveDocument

# --- Step 1: drop the old "_GoBack" bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: locate the "Дополнительные требования(...)" paragraph --------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Дополнительные требовани*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $reqPara  = $d.Paragraphs.Item($targetIndex)
    $nextPara = $d.Paragraphs.Item($targetIndex + 1)

    # Remove the whole following paragraph ("6.1 Адаптивная верстка...",
    # including its paragraph mark) - its text is reinserted below.
    $d.Range($nextPara.Range.Start, $nextPara.Range.End).Delete()

    # Replace the "Дополнительные требования(...)" paragraph's text (but
    # keep its own paragraph mark / numbering formatting) with the new
    # wording.
    $reqPara = $d.Paragraphs.Item($targetIndex)
    $textRange = $d.Range($reqPara.Range.Start, $reqPara.Range.End - 1)
    $textRange.Text = " Адаптивная верстка страниц для мобильного приложения."

    # Re-add the "_GoBack" bookmark, collapsed at the very start of the
    # (now merged) paragraph.
    $reqPara = $d.Paragraphs.Item($targetIndex)
    $bmRange = $d.Range($reqPara.Range.Start, $reqPara.Range.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
